$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 437
$ws.Range("J2").Value = 830
$ws.Range("L2").Value = 830
$ws.Range("N2").Value = -1056

$ws.Range("H132").Value = 6048.25
$ws.Range("I132").Value = 6329.553
$ws.Range("K132").Value = 18988.659
$ws.Range("M132").Value = -16458.659

$ws.Range("H137").Value = 16136203
$ws.Range("I137").Value = 18521100
$ws.Range("K137").Value = 55563300
$ws.Range("M137").Value = -55560750

$ws.Range("H138").Value = 4521.4707
$ws.Range("I138").Value = 4203.9395
$ws.Range("J138").Value = 15000
$ws.Range("K138").Value = 12611.8185
$ws.Range("L138").Value = 45000
$ws.Range("M138").Value = -7471.818500000001
$ws.Range("N138").Value = -55280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1815.4667
$ws.Range("I2").Value = 1269.3334
$ws.Range("K2").Value = 1269.3334
$ws.Range("M2").Value = -1156.3334

$ws.Range("H74").Value = 1990258.2
$ws.Range("I74").Value = 3088424.5
$ws.Range("K74").Value = 3088424.5
$ws.Range("M74").Value = -3087550.5

$ws.Range("H77").Value = 1990258.2
$ws.Range("I77").Value = 3088424.5
$ws.Range("K77").Value = 15442122.5
$ws.Range("M77").Value = -15437754.5

$ws.Range("H92").Value = 36212.75
$ws.Range("J92").Value = 36212.75
$ws.Range("L92").Value = 36212.75
$ws.Range("N92").Value = -41204.75

$ws.Range("H112").Value = 25000
$ws.Range("J112").Value = 25000
$ws.Range("L112").Value = 25000
$ws.Range("N112").Value = -27954

$ws.Range("H116").Value = 1815.4667
$ws.Range("I116").Value = 1269.3334
$ws.Range("K116").Value = 1269.3334
$ws.Range("M116").Value = 1024.6666

$ws.Range("H122").Value = 2304.7693
$ws.Range("I122").Value = 2131.6365
$ws.Range("J122").Value = 3257
$ws.Range("K122").Value = 6394.9095
$ws.Range("L122").Value = 9771
$ws.Range("M122").Value = -3944.9095
$ws.Range("N122").Value = -14671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1815.4667
$ws.Range("I3").Value = 1269.3334
$ws.Range("K3").Value = 1269.3334
$ws.Range("M3").Value = -1155.3334

$ws.Range("H20").Value = 50027.773
$ws.Range("I20").Value = 64228.59
$ws.Range("J20").Value = 1745
$ws.Range("K20").Value = 64228.59
$ws.Range("L20").Value = 1745
$ws.Range("M20").Value = -63981.59
$ws.Range("N20").Value = -2239

$ws.Range("H134").Value = 2876777.5
$ws.Range("I134").Value = 1514.566
$ws.Range("J134").Value = 33354564
$ws.Range("K134").Value = 4543.698
$ws.Range("L134").Value = 100063692
$ws.Range("M134").Value = -2008.698
$ws.Range("N134").Value = -100068762

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 310.25
$ws.Range("I2").Value = 310.25
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 310.25
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -197.25
$ws.Range("N2").ClearContents()

$ws.Range("H22").Value = 1125.4324
$ws.Range("I22").Value = 505.4074
$ws.Range("J22").Value = 2799.5
$ws.Range("K22").Value = 505.4074
$ws.Range("L22").Value = 2799.5
$ws.Range("M22").Value = -155.4074
$ws.Range("N22").Value = -3499.5

$ws.Range("H31").Value = 2144502.5
$ws.Range("I31").Value = 4279005
$ws.Range("K31").Value = 4279005
$ws.Range("M31").Value = -4278710

$ws.Range("H34").Value = 2144502.5
$ws.Range("I34").Value = 4279005
$ws.Range("K34").Value = 4279005
$ws.Range("M34").Value = -4278803

$ws.Range("H50").Value = 23799.75
$ws.Range("J50").Value = 23799.75
$ws.Range("L50").Value = 23799.75
$ws.Range("N50").Value = -25049.75

$ws.Range("H51").Value = 28000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 28000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 28000
$ws.Range("N51").Value = -29472
$ws.Range("M51").ClearContents()

$ws.Range("H60").Value = 21249.5
$ws.Range("I60").Value = 7000
$ws.Range("J60").Value = 25999.334
$ws.Range("K60").Value = 7000
$ws.Range("L60").Value = 25999.334
$ws.Range("M60").Value = -6489
$ws.Range("N60").Value = -27021.334

$ws.Range("H61").Value = 28000
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 28000
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 28000
$ws.Range("N61").Value = -28696
$ws.Range("M61").ClearContents()

$ws.Range("H132").Value = 1917.7179
$ws.Range("I132").Value = 1690.4
$ws.Range("K132").Value = 5071.200000000001
$ws.Range("M132").Value = -2541.200000000001

$ws.Range("H134").Value = 1889.5927
$ws.Range("I134").Value = 1640.76
$ws.Range("K134").Value = 4922.28
$ws.Range("M134").Value = -2387.28

$ws.Range("H141").Value = 73651.625
$ws.Range("J141").Value = 73651.625
$ws.Range("L141").Value = 73651.625
$ws.Range("N141").Value = -84011.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1553190.8
$ws.Range("J5").Value = 1552016.6
$ws.Range("L5").Value = 4656049.800000001
$ws.Range("N5").Value = -4656273.800000001

$ws.Range("H12").Value = 642.5714
$ws.Range("I12").Value = 968.5
$ws.Range("J12").Value = 442
$ws.Range("K12").Value = 2905.5
$ws.Range("L12").Value = 1326
$ws.Range("M12").Value = -2732.5
$ws.Range("N12").Value = -1672

$ws.Range("H64").Value = 6994
$ws.Range("I64").Value = 6994
$ws.Range("K64").Value = 20982
$ws.Range("M64").Value = -20712

$ws.Range("H67").Value = 6994
$ws.Range("I67").Value = 6994
$ws.Range("K67").Value = 20982
$ws.Range("M67").Value = -20046

$ws.Range("H113").Value = 1269.1052
$ws.Range("J113").Value = 1311.8334
$ws.Range("L113").Value = 3935.5002
$ws.Range("N113").Value = -8275.5002

$ws.Range("H114").Value = 2840.5217
$ws.Range("I114").Value = 1325.7142
$ws.Range("J114").Value = 3503.25
$ws.Range("K114").Value = 3977.1426
$ws.Range("L114").Value = 10509.75
$ws.Range("M114").Value = -723.1425999999997
$ws.Range("N114").Value = -17017.75

$ws.Range("H135").Value = 1553190.8
$ws.Range("J135").Value = 1552016.6
$ws.Range("L135").Value = 13968149.4
$ws.Range("N135").Value = -13973219.4

$ws.Range("H137").Value = 3457.2778
$ws.Range("I137").Value = 2304.8333
$ws.Range("J137").Value = 4033.5
$ws.Range("K137").Value = 6914.499899999999
$ws.Range("L137").Value = 12100.5
$ws.Range("M137").Value = -1814.499899999999
$ws.Range("N137").Value = -22300.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2019.0769
$ws.Range("I102").Value = 2019.0769
$ws.Range("K102").Value = 2019.0769
$ws.Range("M102").Value = -397.0769

$ws.Range("H113").Value = 2257.5715
$ws.Range("I113").Value = 2269.6924
$ws.Range("K113").Value = 2269.6924
$ws.Range("M113").Value = -99.69239999999991

$ws.Range("H132").Value = 30368.846
$ws.Range("I132").Value = 30816.25
$ws.Range("J132").Value = 25000
$ws.Range("K132").Value = 92448.75
$ws.Range("L132").Value = 75000
$ws.Range("M132").Value = -89918.75
$ws.Range("N132").Value = -80060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 11053.556
$ws.Range("I61").Value = 11810.25
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 11810.25
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -11608.25
$ws.Range("N61").Value = -5404

$ws.Range("H68").Value = 2000
$ws.Range("I68").Value = 2000
$ws.Range("K68").Value = 2000
$ws.Range("M68").Value = -1251

$ws.Range("H71").Value = 2000
$ws.Range("I71").Value = 2000
$ws.Range("K71").Value = 10000
$ws.Range("M71").Value = -6256

$ws.Range("H93").Value = 1737.8529
$ws.Range("I93").Value = 1359.5385
$ws.Range("K93").Value = 1359.5385
$ws.Range("M93").Value = -111.5385000000001

$ws.Range("H110").Value = 34999
$ws.Range("J110").Value = 34999
$ws.Range("L110").Value = 34999
$ws.Range("N110").Value = -43179

$ws.Range("H113").Value = 11053.556
$ws.Range("I113").Value = 11810.25
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 11810.25
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -9640.25
$ws.Range("N113").Value = -9340

$ws.Range("H136").Value = 8933074
$ws.Range("I136").Value = 5435722.5
$ws.Range("J136").Value = 25020888
$ws.Range("K136").Value = 16307167.5
$ws.Range("L136").Value = 75062664
$ws.Range("M136").Value = -16304617.5
$ws.Range("N136").Value = -75067764

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1532.5454
$ws.Range("I81").Value = 1801
$ws.Range("J81").Value = 816.6667
$ws.Range("K81").Value = 3602
$ws.Range("L81").Value = 1633.3334
$ws.Range("M81").Value = -2541
$ws.Range("N81").Value = -3755.3334

$ws.Range("H84").Value = 1532.5454
$ws.Range("I84").Value = 1801
$ws.Range("J84").Value = 816.6667
$ws.Range("K84").Value = 18010
$ws.Range("L84").Value = 8166.666999999999
$ws.Range("M84").Value = -12706
$ws.Range("N84").Value = -18774.667

$ws.Range("H132").Value = 6411892
$ws.Range("I132").Value = 7577441
$ws.Range("K132").Value = 22732323
$ws.Range("M132").Value = -22729793

$ws.Range("H136").Value = 28261628
$ws.Range("I136").Value = 29814526
$ws.Range("K136").Value = 89443578
$ws.Range("M136").Value = -89441028
